$wb = $excel.ActiveWorkbook

# 1. Rename the "Contingencies" sheet to "Contingencies_Breakers".
#    (This updates the sheet tab name, the defined name range reference,
#    and any formulas pointing at the sheet automatically.)
$contSheet = $wb.Worksheets.Item("Contingencies")
$contSheet.Name = "Contingencies_Breakers"

# 2. Delete row 6 (the "Base_Case" placeholder row) on the Contingencies
#    sheet, shifting the "TEST Cont" / "TEST Cont NC" rows up by one.
$contSheet.Rows.Item(6).Delete()

# 3. Update the active selection / active sheet so the Contingencies
#    sheet is active with C6 selected.
$contSheet.Activate()
$contSheet.Range("C6").Select()

# 4. Update the Loadflow_Settings sheet's selection (it's no longer the
#    active tab) to C31.
$lfSheet = $wb.Worksheets.Item("Loadflow_Settings")
$lfSheet.Range("C31").Select()

# 5. Re-activate the Contingencies_Breakers sheet so it ends up as the
#    workbook's active tab.
$contSheet.Activate()
$contSheet.Range("C6").Select()
